$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = '协鑫集成'
$ws.Range("B2").Value = '协鑫集成'
$ws.Range("C2").Value = '协鑫集成'
$ws.Range("A3").Value = '巨力索具'
$ws.Range("B3").Value = '巨力索具'
$ws.Range("C3").Value = '巨力索具'
$ws.Range("A4").Value = '航天发展'
$ws.Range("B4").Value = '平潭发展'
$ws.Range("C4").Value = '平潭发展'
$ws.Range("A5").Value = '平潭发展'
$ws.Range("B5").Value = '白银有色'
$ws.Range("C5").Value = '利欧股份'
$ws.Range("A6").Value = '浙文互联'
$ws.Range("B6").Value = '利欧股份'
$ws.Range("C6").Value = '神剑股份'
$ws.Range("A7").Value = '利欧股份'
$ws.Range("B7").Value = '航天发展'
$ws.Range("C7").Value = '航天发展'
$ws.Range("A8").Value = '白银有色'
$ws.Range("B8").Value = '浙文互联'
$ws.Range("C8").Value = '白银有色'
$ws.Range("A9").Value = '蓝色光标'
$ws.Range("B9").Value = 'TCL中环'
$ws.Range("C9").Value = '锋龙股份'
$ws.Range("A10").Value = '神剑股份'
$ws.Range("B10").Value = '神剑股份'
$ws.Range("C10").Value = '嘉美包装'
$ws.Range("A11").Value = '金风科技'
$ws.Range("B11").Value = '贵州茅台'
$ws.Range("C11").Value = '浙文互联'
$ws.Range("A12").Value = '三江购物'
$ws.Range("B12").Value = '拓日新能'
$ws.Range("C12").Value = '中超控股'
$ws.Range("A13").Value = '湖南白银'
$ws.Range("B13").Value = '湖南白银'
$ws.Range("C13").Value = '金风科技'
$ws.Range("A14").Value = '茂业商业'
$ws.Range("B14").Value = '茂业商业'
$ws.Range("C14").Value = '湖南白银'
$ws.Range("A15").Value = '网宿科技'
$ws.Range("B15").Value = '金风科技'
$ws.Range("C15").Value = '横店影视'
$ws.Range("A16").Value = '贵州茅台'
$ws.Range("B16").Value = '蓝色光标'
$ws.Range("C16").Value = '东百集团'
$ws.Range("A17").Value = 'TCL中环'
$ws.Range("B17").Value = '海峡创新'
$ws.Range("C17").Value = '茂业商业'
$ws.Range("A18").Value = '海峡创新'
$ws.Range("B18").Value = '华林证券'
$ws.Range("C18").Value = '三江购物'
$ws.Range("A19").Value = '中超控股'
$ws.Range("B19").Value = '横店影视'
$ws.Range("C19").Value = '蓝色光标'
$ws.Range("A20").Value = '锋龙股份'
$ws.Range("B20").Value = '天地在线'
$ws.Range("C20").Value = '海峡创新'
$ws.Range("A21").Value = '新 华 都'
$ws.Range("B21").Value = '遥望科技'
$ws.Range("C21").Value = '杭州解百'
